$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Add the new "Solution" column (D) header + the two existing "solved"
#    bugs' solution text. These are written FIRST (while the rows still
#    carry their original row-level style) so the new D cells inherit the
#    same style that the rest of that row already uses (style 7), and so
#    the shared-string table picks up the new strings in the same order
#    they appear in the authoritative diff (21=Solution, 22=Added a list...,
#    23=Use ClickOnce to deploy, 24=Add ability to deploy ClickOnce...).
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "Solution"
$ws.Range("D12").Value = "Added a list of pre-defined shortcut keys and attempt to register each one "
$ws.Range("D7").Value = "Use ClickOnce to deploy"

# ---------------------------------------------------------------------------
# 2) Insert the new bug row (priority 2) right before the old row 18,
#    pushing "If Start menu is opened..." down to row 19.
# ---------------------------------------------------------------------------
$ws.Rows.Item(18).Insert()
$ws.Range("A18").Value = 2
$ws.Range("B18").Value = "Add ability to deploy ClickOnce without hardcoding publisher name and product name"

# ---------------------------------------------------------------------------
# 3) Mark bug rows 7 and 12 as "solved" (strikethrough + centered priority,
#    strikethrough description) to match rows 2-6/11/13 styling, and row 14
#    likewise gets the strikethrough treatment.
# ---------------------------------------------------------------------------
$ws.Range("A7").Font.Strikethrough = $true
$ws.Range("A7").HorizontalAlignment = -4108
$ws.Range("B7").Font.Strikethrough = $true
$ws.Range("B7").HorizontalAlignment = -4105

$ws.Range("A12").Font.Strikethrough = $true
$ws.Range("A12").HorizontalAlignment = -4108
$ws.Range("B12").Font.Strikethrough = $true
$ws.Range("B12").HorizontalAlignment = -4105
$ws.Range("C12").Clear()

$ws.Range("A14").Font.Strikethrough = $true
$ws.Range("A14").HorizontalAlignment = -4108
$ws.Range("B14").Font.Strikethrough = $true
$ws.Range("B14").HorizontalAlignment = -4105

# ---------------------------------------------------------------------------
# 4) Cosmetics: size the new Solution column, fix the sort range to cover
#    the new column/row, and restore the selection to B9.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 68.7109375

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A19"))
$ws.Sort.SetRange($ws.Range("A2:D19"))
$ws.Sort.Header = -4142

$ws.Range("B9").Select()
